# Template cleanup: trim the generic "pitch deck" scaffolding down to just
# the single slide that the automatic results generator actually reuses
# (the "Results" slide, original sldId 447), and drop its leftover slide
# number placeholder so the generated decks don't show a stray "1" footer.

$p = $ppt.ActivePresentation

# Remove every slide except the one with sldId 447 ("Results"). Iterate
# back-to-front so deleting a slide never shifts the index of one we
# still need to inspect.
for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $slide = $p.Slides.Item($i)
    if ($slide.SlideID -ne 447) {
        $slide.Delete()
    }
}

# The lone remaining slide is now the only item in the collection.
$s = $p.Slides.Item(1)

# Drop the "Slide Number Placeholder 2" shape. Toggling the header/footer
# visibility off (rather than calling Shape.Delete on the placeholder)
# removes it from the slide's shape tree instead of leaving an inherited
# placeholder ghost behind.
$s.HeadersFooters.SlideNumber.Visible = [Microsoft.Office.Core.MsoTriState]::msoFalse
